$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (columns D and E inserted: ownTeam, oppTeam)
$ws.Range("A1").Value = "venue"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "result"
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"
$ws.Range("F1").Value = "batsman"
$ws.Range("G1").Value = "totalRuns"
$ws.Range("H1").Value = "totalBalls"
$ws.Range("I1").Value = "total4s"
$ws.Range("J1").Value = "total6s"
$ws.Range("K1").Value = "sr"

$data = @(
    @(" Abu Dhabi", " October 16 2020", "Mumbai won by 8 wickets (with 19 balls remaining)", "Kolkata Knight Riders", "Mumbai Indians", "Eoin Morgan (c)", "39", "29", "2", "2", "134.48"),
    @(" Abu Dhabi", " October 18 2020", "Match tied (KKR won the one-over eliminator)", "Kolkata Knight Riders", "Sunrisers Hyderabad", "Eoin Morgan (c)", "34", "23", "3", "1", "147.82"),
    @(" Dubai (DSC)", " November 01 2020", "KKR won by 60 runs", "Kolkata Knight Riders", "Rajasthan Royals", "Eoin Morgan (c)", "68", "35", "5", "6", "194.28"),
    @(" Dubai (DSC)", " October 29 2020", "Super Kings won by 6 wickets", "Kolkata Knight Riders", "Chennai Super Kings", "Eoin Morgan (c)", "15", "12", "2", "0", "125.00"),
    @(" Sharjah", " October 26 2020", "Kings XI won by 8 wickets (with 7 balls remaining)", "Kolkata Knight Riders", "Kings XI Punjab", "Eoin Morgan (c)", "40", "25", "5", "2", "160.00"),
    @(" Abu Dhabi", " October 21 2020", "RCB won by 8 wickets (with 39 balls remaining)", "Kolkata Knight Riders", "Royal Challengers Bangalore", "Eoin Morgan (c)", "30", "34", "3", "1", "88.23"),
    @(" Abu Dhabi", " October 24 2020", "KKR won by 59 runs", "Kolkata Knight Riders", "Delhi Capitals", "Eoin Morgan (c)", "17", "9", "2", "1", "188.88")
)

# Ensure numeric-looking columns (G:K) stay stored as text, matching source data
$ws.Range("G2:K8").NumberFormat = "@"

$row = 2
foreach ($r in $data) {
    $col = 1
    foreach ($val in $r) {
        $ws.Cells.Item($row, $col).Value = $val
        $col++
    }
    $row++
}
